$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,@("D2", "36.479.51")
    ,@("E2", "  +0.13%  ")
    ,@("D3", "1.941.43")
    ,@("E3", "  -1.35%  ")
    ,@("E4", "  +0.01%  ")
    ,@("D5", "243.73")
    ,@("E5", "  -0.37%  ")
    ,@("E6", "  -0.86%  ")
    ,@("E7", "  +0.01%  ")
    ,@("D8", "57.35")
    ,@("E8", "  -1.57%  ")
    ,@("D9", "0.362")
    ,@("E9", "  -2.60%  ")
    ,@("D10", "0.0849")
    ,@("E10", "  +1.28%  ")
    ,@("E11", "  -0.78%  ")
    ,@("D12", "2.228.18")
    ,@("E12", "  -1.19%  ")
    ,@("D13", "21.28")
    ,@("E13", "  -4.42%  ")
    ,@("D14", "0.810")
    ,@("E14", "  -2.72%  ")
    ,@("D15", "13.44")
    ,@("E15", "  -0.96%  ")
    ,@("E16", "  -3.49%  ")
    ,@("D17", "1.940.00")
    ,@("E17", "  -1.91%  ")
    ,@("D18", "36.438.58")
    ,@("E18", "  +0.32%  ")
    ,@("D19", "69.37")
    ,@("E19", "  -2.71%  ")
    ,@("D20", "0.0$([char]8323)0864")
    ,@("E20", "  -2.13%  ")
    ,@("D21", "228.13")
    ,@("E21", "  -0.71%  ")
    ,@("E22", "  -2.92%  ")
    ,@("E23", "  +0.02%  ")
    ,@("D24", "2.36")
    ,@("E24", "  -6.33%  ")
    ,@("E25", "  +1.27%  ")
    ,@("E26", "  -4.58%  ")
    ,@("D27", "161.13")
    ,@("E27", "  -2.49%  ")
    ,@("E28", "  +8.40%  ")
    ,@("D29", "19.17")
    ,@("E29", "  -3.90%  ")
    ,@("E30", "  -0.86%  ")
    ,@("E31", "  -5.15%  ")
    ,@("D32", "4.58")
    ,@("E32", "  -3.78%  ")
    ,@("D33", "0.0616")
    ,@("E33", "  -3.83%  ")
    ,@("E34", "  -3.34%  ")
    ,@("E35", "  +4.28%  ")
    ,@("E36", "  +0.06%  ")
    ,@("E37", "  -1.74%  ")
    ,@("E38", "  +0.42%  ")
    ,@("D39", "3.19")
    ,@("E39", "  +8.68%  ")
    ,@("D40", "0.0983")
    ,@("E40", "  +1.81%  ")
    ,@("E41", "  +0.63%  ")
    ,@("E42", "  -2.87%  ")
    ,@("D43", "0.0208")
    ,@("E43", "  -1.20%  ")
    ,@("D44", "15.99")
    ,@("E44", "  +1.29%  ")
    ,@("D45", "1.342.15")
    ,@("E45", "  -0.69%  ")
    ,@("E46", "  -3.22%  ")
    ,@("D47", "86.58")
    ,@("E47", "  -2.89%  ")
    ,@("D48", "7.20")
    ,@("E48", "  -1.14%  ")
    ,@("D49", "2.81")
    ,@("E49", "  -0.51%  ")
    ,@("D50", "2.120.26")
    ,@("E50", "  -1.12%  ")
    ,@("D51", "43.23")
)

foreach ($change in $changes) {
    $targetCell = $ws.Range($change[0])
    $targetCell.NumberFormat = "@"
    $targetCell.Value = $change[1]
}
